# semana 06 de 2026
# Update the poisson.xlsx "Esperado/Observado/valor p" table:
#  - insert a new event row (610) right before the "Parotiditis" (620) row,
#    which shifts the rows below it down by one
#  - refresh the Esperado / Observado / valor p figures for this week

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "610" row above row 22 (old row 22 -> new row 23, etc.) ---
$ws.Rows.Item(22).Insert()

# New row only carries an event code (text, like the rest of column A) and an
# Observado value of 0 - Esperado/nom_eve/valor p are left blank for this event.
# (Route the literal through a formula + paste-values so it lands as text,
# matching the rest of column A, instead of Excel auto-coercing "610" to a number.)
$ws.Range("A22").Formula = '="610"'
$ws.Range("A22").Copy()
$ws.Range("A22").PasteSpecial(-4163)
$ws.Range("D22").Value = 0

# --- Refresh Esperado (C) / Observado (D) / valor p (E) figures ---
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 0.27

$ws.Range("D3").Value = 1

$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 11
$ws.Range("E4").Value = 0.01

$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 31

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.27

$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1

$ws.Range("C8").Value = 46
$ws.Range("D8").Value = 39
$ws.Range("E8").Value = 0.04

$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 0.27

$ws.Range("C11").Value = 67

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = 0.1

$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 0.37

$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 0.09

$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 1

$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 0.09

$ws.Range("D20").Value = 0

$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0

$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 0.37

$ws.Range("C24").Value = 0
$ws.Range("E24").Value = 0

$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0.18

$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 0.07

$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 0.06

$ws.Range("C28").Value = 8
$ws.Range("D28").Value = 8
$ws.Range("E28").Value = 0.14
